$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new MAC-address groups (10030 / 10031), 5 rows each, continuing the
# existing device_id sequence from 3000165.
$machineIds = @(10030, 10030, 10030, 10030, 10030, 10031, 10031, 10031, 10031, 10031)
$deviceIdStart = 3000166

$startRow = 147
for ($i = 0; $i -lt $machineIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $machineIds[$i]
    $ws.Cells.Item($row, 3).Value = $deviceIdStart + $i
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
}

$ws.Range("A148").Select()

# Mirror the author's scrolled viewport (best-effort; harmless if unsupported).
try {
    $excel.ActiveWindow.ScrollRow = 142
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
